$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1735
$ws.Range("F5").Value = 64
$ws.Range("F6").Value = 701
$ws.Range("F8").Value = 200
$ws.Range("F11").Value = 31
$ws.Range("F12").Value = 562
$ws.Range("F13").Value = 485
$ws.Range("F14").Value = 132
$ws.Range("F19").Value = 518
$ws.Range("F23").Value = 189
$ws.Range("F25").Value = 133
$ws.Range("F28").Value = 51
$ws.Range("F29").Value = 158
$ws.Range("F32").Value = 34
$ws.Range("F34").Value = 33

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 343
$ws.Range("F23").Value = 263

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2352
$ws.Range("F9").Value = 1192

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 2352
$ws.Range("F7").Value = 1735
$ws.Range("F10").Value = 1192
$ws.Range("F14").Value = 64
$ws.Range("F15").Value = 701
$ws.Range("F18").Value = 200
$ws.Range("F20").Value = 31
$ws.Range("F21").Value = 562
$ws.Range("F22").Value = 485
$ws.Range("F23").Value = 132
$ws.Range("F28").Value = 518
$ws.Range("F32").Value = 189
$ws.Range("F33").Value = 133
$ws.Range("F38").Value = 51
$ws.Range("F39").Value = 158
$ws.Range("F44").Value = 263
$ws.Range("F48").Value = 34
$ws.Range("F50").Value = 33
